# Updated cryptos list - applies Price (D) and Volume(1h) (E) changes per row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to hold a literal text value (matches the source
    # workbook which stores these as inline/shared strings), even when
    # the text happens to look like a number (e.g. "213.55").
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "27.785.51"
Set-TextCell $ws.Range("E2") "  +1.39%  "

Set-TextCell $ws.Range("D3") "1.650.54"
Set-TextCell $ws.Range("E3") "  -0.23%  "

Set-TextCell $ws.Range("E4") "  +0.07%  "

Set-TextCell $ws.Range("D5") "213.55"
Set-TextCell $ws.Range("E5") "  +0.14%  "

Set-TextCell $ws.Range("E6") "  -0.71%  "

Set-TextCell $ws.Range("E7") "  +0.10%  "

Set-TextCell $ws.Range("D8") "23.20"
Set-TextCell $ws.Range("E8") "  -1.01%  "

Set-TextCell $ws.Range("E9") "  -0.35%  "

Set-TextCell $ws.Range("D10") "0.0616"
Set-TextCell $ws.Range("E10") "  +0.44%  "

Set-TextCell $ws.Range("E11") "  -1.64%  "

Set-TextCell $ws.Range("D12") "1.884.20"
Set-TextCell $ws.Range("E12") "  -0.26%  "

Set-TextCell $ws.Range("D13") "1.652.05"
Set-TextCell $ws.Range("E13") "  -0.12%  "

Set-TextCell $ws.Range("E14") "  -0.43%  "

Set-TextCell $ws.Range("E15") "  -0.31%  "

Set-TextCell $ws.Range("D16") "64.57"
Set-TextCell $ws.Range("E16") "  -1.42%  "

Set-TextCell $ws.Range("D17") "27.768.98"
Set-TextCell $ws.Range("E17") "  +1.37%  "

Set-TextCell $ws.Range("D18") "234.22"
Set-TextCell $ws.Range("E18") "  +2.30%  "

Set-TextCell $ws.Range("D19") "7.72"
Set-TextCell $ws.Range("E19") "  +4.07%  "

Set-TextCell $ws.Range("E20") "  -0.08%  "

Set-TextCell $ws.Range("E21") "  +0.10%  "

Set-TextCell $ws.Range("E22") "  -0.37%  "

Set-TextCell $ws.Range("D23") "10.16"
Set-TextCell $ws.Range("E23") "  +8.25%  "

Set-TextCell $ws.Range("E24") "  -3.74%  "

Set-TextCell $ws.Range("D25") "150.63"
Set-TextCell $ws.Range("E25") "  +2.44%  "

Set-TextCell $ws.Range("E26") "  -0.94%  "

Set-TextCell $ws.Range("E27") "  -1.73%  "

Set-TextCell $ws.Range("E28") "  +0.15%  "

Set-TextCell $ws.Range("E29") "  +0.10%  "

Set-TextCell $ws.Range("E30") "  +0.31%  "

Set-TextCell $ws.Range("E31") "  -1.06%  "

Set-TextCell $ws.Range("E32") "  +0.55%  "

Set-TextCell $ws.Range("E33") "  +1.60%  "

Set-TextCell $ws.Range("D34") "1.443.63"
Set-TextCell $ws.Range("E34") "  +1.56%  "

Set-TextCell $ws.Range("D35") "1.59"
Set-TextCell $ws.Range("E35") "  +2.25%  "

Set-TextCell $ws.Range("E36") "  -1.06%  "

Set-TextCell $ws.Range("D37") "0.571"
Set-TextCell $ws.Range("E37") "  +0.70%  "

Set-TextCell $ws.Range("D38") "0.888"
Set-TextCell $ws.Range("E38") "  -1.89%  "

Set-TextCell $ws.Range("E39") "  -0.60%  "

Set-TextCell $ws.Range("D40") "0.879"
Set-TextCell $ws.Range("E40") "  +11.39%  "

Set-TextCell $ws.Range("E41") "  -0.67%  "

Set-TextCell $ws.Range("E42") "  +0.11%  "

Set-TextCell $ws.Range("D43") "5.62"
Set-TextCell $ws.Range("E43") "  +1.50%  "

Set-TextCell $ws.Range("D44") "66.69"
Set-TextCell $ws.Range("E44") "  +2.58%  "

Set-TextCell $ws.Range("D45") "2.47"
Set-TextCell $ws.Range("E45") "  -1.21%  "

Set-TextCell $ws.Range("E46") "  +2.23%  "

Set-TextCell $ws.Range("D47") "1.792.46"
Set-TextCell $ws.Range("E47") "  -0.32%  "

Set-TextCell $ws.Range("D48") "1.74"
Set-TextCell $ws.Range("E48") "  +4.32%  "

Set-TextCell $ws.Range("D49") "86.50"
Set-TextCell $ws.Range("E49") "  -1.67%  "

Set-TextCell $ws.Range("D50") "0.0₆0106"
Set-TextCell $ws.Range("E50") "  +1.34%  "

Set-TextCell $ws.Range("D51") "0.0999"
Set-TextCell $ws.Range("E51") "  -1.18%  "

